$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44400
$ws.Range("J2").Value = 70
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 15000
$ws.Range("P2").Value = 500

# Row 3
$ws.Range("D3").Value = 44162
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("P3").Value = 250

# Row 4
$ws.Range("D4").Value = 44176
$ws.Range("J4").Value = 80

# Row 5
$ws.Range("D5").Value = 44446
$ws.Range("H5").Value = 'Madrigal'
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 16000
$ws.Range("L5").Value = 16000
$ws.Range("M5").Value = 16000
$ws.Range("N5").Value = '$/caja 40 unidades'
$ws.Range("P5").Value = 400
$ws.Range("Q5").Value = 40

# Row 6
$ws.Range("D6").Value = 44418
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 16000
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 16000
$ws.Range("P6").Value = 533

# Row 7
$ws.Range("D7").Value = 44421
$ws.Range("K7").Value = 16500
$ws.Range("L7").Value = 16500
$ws.Range("M7").Value = 16500
$ws.Range("P7").Value = 550

# Row 8
$ws.Range("D8").Value = 44161
$ws.Range("H8").Value = 'Madrigal'
$ws.Range("J8").Value = 30
$ws.Range("K8").Value = 11000
$ws.Range("L8").Value = 11000
$ws.Range("M8").Value = 11000
$ws.Range("N8").Value = '$/caja 40 unidades'
$ws.Range("P8").Value = 275
$ws.Range("Q8").Value = 40

# Row 9
$ws.Range("D9").Value = 44407
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 18000
$ws.Range("L9").Value = 18000
$ws.Range("M9").Value = 18000
$ws.Range("P9").Value = 600

# Row 11
$ws.Range("D11").Value = 44386
$ws.Range("H11").Value = 'Española'
$ws.Range("J11").Value = 30
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 15000
$ws.Range("N11").Value = '$/caja 30 unidades'
$ws.Range("P11").Value = 500
$ws.Range("Q11").Value = 30

# Row 12
$ws.Range("D12").Value = 44390
$ws.Range("H12").Value = 'Española'
$ws.Range("K12").Value = 16000
$ws.Range("L12").Value = 16000
$ws.Range("M12").Value = 16000
$ws.Range("N12").Value = '$/caja 30 unidades'
$ws.Range("P12").Value = 533
$ws.Range("Q12").Value = 30
